$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.65505065700124776
$ws.Range("L1").Value = 0.821840996564285
$ws.Range("BP1").Value = 0.91293480306456509
$ws.Range("B3").Value = 0.88362970599851076
$ws.Range("D3").Value = 0.84931068305072677
$ws.Range("B4").Value = 0.7493031366899916
$ws.Range("F4").Value = 0.6229915625148813
$ws.Range("G4").Value = 0.97386915505288907
$ws.Range("C5").Value = 0.86777631746887851
$ws.Range("BG6").Value = 0.97723102837704812
$ws.Range("E7").Value = 0.8628832385403099
$ws.Range("F7").Value = 0.67694591058610754
$ws.Range("G8").Value = 0.62066486892163597
$ws.Range("I8").Value = 0.74958738177630679
$ws.Range("P8").Value = 0.7358053826940647
$ws.Range("Z9").Value = 0.71515726756553888
$ws.Range("H10").Value = 0.82663699660928192
$ws.Range("I10").Value = 0.61104156441508695
$ws.Range("K10").Value = 0.96996196731853379
$ws.Range("L11").Value = 0.91936563995652143
$ws.Range("K13").Value = 0.90382503265164793
$ws.Range("BF13").Value = 0.63524844502094813
$ws.Range("L14").Value = 0.79238083074983989
$ws.Range("M14").Value = 0.96664944369733552
$ws.Range("N15").Value = 0.98560220020105582
$ws.Range("Q15").Value = 0.99407860905500844
$ws.Range("BG15").Value = 0.84848750318559174
$ws.Range("N16").Value = 0.98514752509992509
$ws.Range("O16").Value = 0.59465931784375381
$ws.Range("AD16").Value = 0.9724266307802284
$ws.Range("F17").Value = 0.71829164516216382
$ws.Range("T18").Value = 0.90510538503210025
$ws.Range("Q19").Value = 0.91256643657661796
$ws.Range("T19").Value = 0.99864019206798527
$ws.Range("U20").Value = 0.98287021745740055
$ws.Range("V20").Value = 0.90340675219949274
$ws.Range("V21").Value = 0.97009961524990773
$ws.Range("S22").Value = 0.72052615874387782
$ws.Range("W22").Value = 0.67396039388308671
$ws.Range("L23").Value = 0.92637289993889693
$ws.Range("X23").Value = 0.92624484576375132
$ws.Range("BP23").Value = 0.62512193717966746
$ws.Range("V24").Value = 0.82198156528312993
$ws.Range("W25").Value = 0.83553256736832227
$ws.Range("X25").Value = 0.75278350282022677
$ws.Range("Z25").Value = 0.71830264785355136
$ws.Range("X26").Value = 0.61338067321737066
$ws.Range("AB26").Value = 0.92991385869319265
$ws.Range("Z27").Value = 0.82666629421558446
$ws.Range("P28").Value = 0.81665368388616377
$ws.Range("AA28").Value = 0.93674745988878372
$ws.Range("AC28").Value = 0.90703407469897335
$ws.Range("AD28").Value = 0.69576619818841468
$ws.Range("U29").Value = 0.94369699618682068
$ws.Range("AA29").Value = 0.95256167502627087
$ws.Range("BO29").Value = 0.93432948546473482
$ws.Range("AC30").Value = 0.76939197714725149
$ws.Range("AE30").Value = 0.67044518150701093
$ws.Range("BE30").Value = 0.61670094988014856
$ws.Range("AC31").Value = 0.95713191333526915
$ws.Range("AF31").Value = 0.85482512322563931
$ws.Range("P32").Value = 0.77488756483333332
$ws.Range("K33").Value = 0.95068096915793676
$ws.Range("AE33").Value = 0.70685902177486648
$ws.Range("AZ33").Value = 0.86224685803940493
$ws.Range("AF34").Value = 0.91203188499513943
$ws.Range("G35").Value = 0.91300789262160476
$ws.Range("AG35").Value = 0.91642549849297539
$ws.Range("AH35").Value = 0.75143796014081243
$ws.Range("AI36").Value = 0.90333546301390344
$ws.Range("AL36").Value = 0.9753084733973767
$ws.Range("R37").Value = 0.71244868183133958
$ws.Range("AI37").Value = 0.7748006117113323
$ws.Range("AJ37").Value = 0.6097169399386333
$ws.Range("AL37").Value = 0.97275964554369077
$ws.Range("AM37").Value = 0.63754044588504966
$ws.Range("AM38").Value = 0.87415062980386971
$ws.Range("AN38").Value = 0.77190971135714292
$ws.Range("W39").Value = 0.90991130060901959
$ws.Range("AN39").Value = 0.70871332140413612
$ws.Range("AO40").Value = 0.98558760570001724
$ws.Range("AP40").Value = 0.9216616590148432
$ws.Range("AP41").Value = 0.89785902334879764
$ws.Range("AQ42").Value = 0.77257587968435448
$ws.Range("V43").Value = 0.78344776494675905
$ws.Range("AO43").Value = 0.91590304821958379
$ws.Range("AP44").Value = 0.74221987481762253
$ws.Range("AS44").Value = 0.90329733184608574
$ws.Range("AK45").Value = 0.62106453353264957
$ws.Range("AQ45").Value = 0.86643733781959398
$ws.Range("AR46").Value = 0.74134856898807966
$ws.Range("AU46").Value = 0.78105582294499099
$ws.Range("AV46").Value = 0.92003526087168197
$ws.Range("AV47").Value = 0.92981099932253453
$ws.Range("AW47").Value = 0.82166080574305389
$ws.Range("AS48").Value = 0.9286553043304776
$ws.Range("AX48").Value = 0.86413816157231038
$ws.Range("AV49").Value = 0.82383075279084128
$ws.Range("AX49").Value = 0.96761695718002483
$ws.Range("AY49").Value = 0.82537714125900585
$ws.Range("J50").Value = 0.75910459453657242
$ws.Range("AZ50").Value = 0.64510317525378402
$ws.Range("AX51").Value = 0.86551730435498619
$ws.Range("AZ51").Value = 0.68708093442682139
$ws.Range("BB52").Value = 0.81186472187875958
$ws.Range("P53").Value = 0.77930113138325385
$ws.Range("AY53").Value = 0.78290166303728559
$ws.Range("AP54").Value = 0.82257124161649631
$ws.Range("BA54").Value = 0.65648749932275263
$ws.Range("BD55").Value = 0.98663203154499945
$ws.Range("BB56").Value = 0.84298581764974734
$ws.Range("BE56").Value = 0.9885081585435983
$ws.Range("BF56").Value = 0.68860015677568054
$ws.Range("BC57").Value = 0.88960759643519793
$ws.Range("BF57").Value = 0.6337693718138635
$ws.Range("Y58").Value = 0.61983811100125341
$ws.Range("BG58").Value = 0.74260005885462776
$ws.Range("BF60").Value = 0.94982049974785032
$ws.Range("BI60").Value = 0.97822365476206063
$ws.Range("AU61").Value = 0.70029037392352911
$ws.Range("BG61").Value = 0.90431575044449941
$ws.Range("BH62").Value = 0.97098480633141437
$ws.Range("BK62").Value = 0.80512521931213832
$ws.Range("BL62").Value = 0.68407717792227374
$ws.Range("BI63").Value = 0.66891984245340685
$ws.Range("BM63").Value = 0.92273136997414551
$ws.Range("E64").Value = 0.64203394489087695
$ws.Range("AJ64").Value = 0.92430781460713263
$ws.Range("BK64").Value = 0.6786128267796161
$ws.Range("BM64").Value = 0.76701510898552416
$ws.Range("BO65").Value = 0.94643037551235865
$ws.Range("AH66").Value = 0.87106246550090105
$ws.Range("BM66").Value = 0.83040729147987913
$ws.Range("BO66").Value = 0.82709113943720747
$ws.Range("BP66").Value = 0.939854575000306
$ws.Range("BO68").Value = 0.96803542607475679
